# Applies the commit's changes to PlayerPerformance_6204.xlsx:
#  1. Remove the (empty) INNING_NUMBER cell from rows 2,3,6,7,8 on "ODI Batting"
#     (these were blank placeholder cells; dropping them leaves a true gap in
#     the row instead of an explicit-but-empty cell).
#  2. Add a new "ODI Batting Extra" worksheet at the end of the workbook with
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#     MAN_OF_MATCH columns of per-match batting extras.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty B cells on "ODI Batting" -----------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B3").ClearContents()
$odiBatting.Range("B6").ClearContents()
$odiBatting.Range("B7").ClearContents()
$odiBatting.Range("B8").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after the last sheet -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Match the bold / bordered / centered header style used by the other sheets
$odiBatting.Range("A1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are text columns; a leading
# apostrophe keeps the digit-looking strings ("4669", "0", "0.52%") from being
# coerced into numbers. BATTING_POSITION is the one genuinely numeric column.
$extra.Range("A2").Value = "'4669"
$extra.Range("B2").Value = 9
$extra.Range("F2").Value = "NO"

$extra.Range("A3").Value = "'4673"
$extra.Range("F3").Value = "NO"

$extra.Range("A4").Value = "'4676"
$extra.Range("B4").Value = 11
$extra.Range("C4").Value = "'0"
$extra.Range("D4").Value = "'0"
$extra.Range("F4").Value = "NO"

$extra.Range("A5").Value = "'4682"
$extra.Range("F5").Value = "NO"

$extra.Range("A6").Value = "'4685"
$extra.Range("F6").Value = "NO"

$extra.Range("A7").Value = "'4687"
$extra.Range("B7").Value = 10
$extra.Range("C7").Value = "'0"
$extra.Range("D7").Value = "'0"
$extra.Range("E7").Value = "'0.52%"
$extra.Range("F7").Value = "NO"

$extra.Range("A8").Value = "'4689"

$extra.Range("A9").Value = "'4697"

# The leading apostrophe marks these as "text stored as number" (quote-prefix)
# formatting. Re-paste plain formatting from an ordinary, never-prefixed text
# cell so the cells keep their Text value but fall back to the default style,
# matching the rest of the sheet.
$wb.Worksheets.Item("Player Info").Range("B2").Copy()
$extra.Range("A2:A9").PasteSpecial(-4122)  # xlPasteFormats
$extra.Range("C4:D4").PasteSpecial(-4122)
$extra.Range("C7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Adding a sheet makes it the active tab; restore the original active sheet
# ("Player Info") so the workbook-level view state is unchanged, matching
# the source workbook (activeTab="0").
$wb.Worksheets.Item(1).Activate()
